$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(2, 7).Value = 7.102474
$ws.Cells.Item(2, 8).Value = 21.307422
$ws.Cells.Item(2, 9).Value = 0.3851819652723766
$ws.Cells.Item(2, 10).Value = 0.3851819652723767
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.1474273333333333
$ws.Cells.Item(2, 14).Value = 0.442282
$ws.Cells.Item(2, 15).Value = 0.1588601259223368
$ws.Cells.Item(2, 16).Value = 0.1588601259223368
$ws.Cells.Item(2, 17).Value = 1.047098801889333
$ws.Cells.Item(2, 18).Value = 9.423889217004
$ws.Cells.Item(2, 19).Value = 0.06119005550618289
$ws.Cells.Item(2, 20).Value = 0.0611900555061829
$ws.Cells.Item(3, 7).Value = 7.102474
$ws.Cells.Item(3, 8).Value = 21.307422
$ws.Cells.Item(3, 9).Value = 0.3851819652723766
$ws.Cells.Item(3, 10).Value = 0.3851819652723767
$ws.Cells.Item(3, 15).Value = 0.4626735347223893
$ws.Cells.Item(3, 16).Value = 0.4626735347223893
$ws.Cells.Item(3, 17).Value = 3.049631876224
$ws.Cells.Item(3, 18).Value = 27.446686886016
$ws.Cells.Item(3, 19).Value = 0.1782135013838871
$ws.Cells.Item(3, 20).Value = 0.1782135013838871
$ws.Cells.Item(4, 7).Value = 7.102474
$ws.Cells.Item(4, 8).Value = 21.307422
$ws.Cells.Item(4, 9).Value = 0.3851819652723766
$ws.Cells.Item(4, 10).Value = 0.3851819652723767
$ws.Cells.Item(4, 13).Value = 0.01780266666666666
$ws.Cells.Item(4, 14).Value = 0.053408
$ws.Cells.Item(4, 15).Value = 0.01918323966442261
$ws.Cells.Item(4, 16).Value = 0.01918323966442261
$ws.Cells.Item(4, 17).Value = 0.1264429771306667
$ws.Cells.Item(4, 18).Value = 1.137986794176
$ws.Cells.Item(4, 19).Value = 0.007389037954233307
$ws.Cells.Item(4, 20).Value = 0.00738903795423331
$ws.Cells.Item(5, 7).Value = 7.102474
$ws.Cells.Item(5, 8).Value = 21.307422
$ws.Cells.Item(5, 9).Value = 0.3851819652723766
$ws.Cells.Item(5, 10).Value = 0.3851819652723767
$ws.Cells.Item(5, 13).Value = 0.3334263333333333
$ws.Cells.Item(5, 14).Value = 1.000279
$ws.Cells.Item(5, 15).Value = 0.3592830996908513
$ws.Cells.Item(5, 16).Value = 0.3592830996908513
$ws.Cells.Item(5, 17).Value = 2.368151863415333
$ws.Cells.Item(5, 18).Value = 21.313366770738
$ws.Cells.Item(5, 19).Value = 0.1383893704280733
$ws.Cells.Item(5, 20).Value = 0.1383893704280733
$ws.Cells.Item(6, 7).Value = 8.299810000000001
$ws.Cells.Item(6, 9).Value = 0.4501159915808667
$ws.Cells.Item(6, 10).Value = 0.4501159915808668
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.1474273333333333
$ws.Cells.Item(6, 14).Value = 0.442282
$ws.Cells.Item(6, 15).Value = 0.1588601259223368
$ws.Cells.Item(6, 16).Value = 0.1588601259223368
$ws.Cells.Item(6, 17).Value = 1.223618855473333
$ws.Cells.Item(6, 18).Value = 11.01256969926
$ws.Cells.Item(6, 19).Value = 0.07150548310219396
$ws.Cells.Item(6, 20).Value = 0.07150548310219397
$ws.Cells.Item(7, 7).Value = 8.299810000000001
$ws.Cells.Item(7, 9).Value = 0.4501159915808667
$ws.Cells.Item(7, 10).Value = 0.4501159915808668
$ws.Cells.Item(7, 15).Value = 0.4626735347223893
$ws.Cells.Item(7, 16).Value = 0.4626735347223893
$ws.Cells.Item(7, 17).Value = 3.56373921856
$ws.Cells.Item(7, 19).Value = 0.2082567568597928
$ws.Cells.Item(7, 20).Value = 0.2082567568597929
$ws.Cells.Item(8, 7).Value = 8.299810000000001
$ws.Cells.Item(8, 9).Value = 0.4501159915808667
$ws.Cells.Item(8, 10).Value = 0.4501159915808668
$ws.Cells.Item(8, 13).Value = 0.01780266666666666
$ws.Cells.Item(8, 14).Value = 0.053408
$ws.Cells.Item(8, 15).Value = 0.01918323966442261
$ws.Cells.Item(8, 16).Value = 0.01918323966442261
$ws.Cells.Item(8, 17).Value = 0.1477587508266667
$ws.Cells.Item(8, 18).Value = 1.32982875744
$ws.Cells.Item(8, 19).Value = 0.008634682943284995
$ws.Cells.Item(8, 20).Value = 0.008634682943284998
$ws.Cells.Item(9, 7).Value = 8.299810000000001
$ws.Cells.Item(9, 9).Value = 0.4501159915808667
$ws.Cells.Item(9, 10).Value = 0.4501159915808668
$ws.Cells.Item(9, 13).Value = 0.3334263333333333
$ws.Cells.Item(9, 14).Value = 1.000279
$ws.Cells.Item(9, 15).Value = 0.3592830996908513
$ws.Cells.Item(9, 16).Value = 0.3592830996908513
$ws.Cells.Item(9, 17).Value = 2.767375215663333
$ws.Cells.Item(9, 18).Value = 24.90637694097
$ws.Cells.Item(9, 19).Value = 0.1617190686755949
$ws.Cells.Item(9, 20).Value = 0.1617190686755949
$ws.Cells.Item(10, 7).Value = 0.2555593333333333
$ws.Cells.Item(10, 8).Value = 0.766678
$ws.Cells.Item(10, 9).Value = 0.01385951518541732
$ws.Cells.Item(10, 10).Value = 0.01385951518541733
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.1474273333333333
$ws.Cells.Item(10, 14).Value = 0.442282
$ws.Cells.Item(10, 15).Value = 0.1588601259223368
$ws.Cells.Item(10, 16).Value = 0.1588601259223368
$ws.Cells.Item(10, 17).Value = 0.03767643102177777
$ws.Cells.Item(10, 18).Value = 0.339087879196
$ws.Cells.Item(10, 19).Value = 0.002201724327577934
$ws.Cells.Item(10, 20).Value = 0.002201724327577935
$ws.Cells.Item(11, 7).Value = 0.2555593333333333
$ws.Cells.Item(11, 8).Value = 0.766678
$ws.Cells.Item(11, 9).Value = 0.01385951518541732
$ws.Cells.Item(11, 10).Value = 0.01385951518541733
$ws.Cells.Item(11, 15).Value = 0.4626735347223893
$ws.Cells.Item(11, 16).Value = 0.4626735347223893
$ws.Cells.Item(11, 17).Value = 0.1097310443093333
$ws.Cells.Item(11, 18).Value = 0.9875793987839999
$ws.Cells.Item(11, 19).Value = 0.006412430880375663
$ws.Cells.Item(11, 20).Value = 0.006412430880375666
$ws.Cells.Item(12, 7).Value = 0.2555593333333333
$ws.Cells.Item(12, 8).Value = 0.766678
$ws.Cells.Item(12, 9).Value = 0.01385951518541732
$ws.Cells.Item(12, 10).Value = 0.01385951518541733
$ws.Cells.Item(12, 13).Value = 0.01780266666666666
$ws.Cells.Item(12, 14).Value = 0.053408
$ws.Cells.Item(12, 15).Value = 0.01918323966442261
$ws.Cells.Item(12, 16).Value = 0.01918323966442261
$ws.Cells.Item(12, 17).Value = 0.004549637624888888
$ws.Cells.Item(12, 18).Value = 0.04094673862399999
$ws.Cells.Item(12, 19).Value = 0.000265870401434565
$ws.Cells.Item(12, 20).Value = 0.0002658704014345652
$ws.Cells.Item(13, 7).Value = 0.2555593333333333
$ws.Cells.Item(13, 8).Value = 0.766678
$ws.Cells.Item(13, 9).Value = 0.01385951518541732
$ws.Cells.Item(13, 10).Value = 0.01385951518541733
$ws.Cells.Item(13, 13).Value = 0.3334263333333333
$ws.Cells.Item(13, 14).Value = 1.000279
$ws.Cells.Item(13, 15).Value = 0.3592830996908513
$ws.Cells.Item(13, 16).Value = 0.3592830996908513
$ws.Cells.Item(13, 17).Value = 0.08521021146244444
$ws.Cells.Item(13, 18).Value = 0.766891903162
$ws.Cells.Item(13, 19).Value = 0.004979489576029159
$ws.Cells.Item(13, 20).Value = 0.00497948957602916
$ws.Cells.Item(14, 7).Value = 1.392645666666667
$ws.Cells.Item(14, 8).Value = 4.177937
$ws.Cells.Item(14, 9).Value = 0.07552607652132563
$ws.Cells.Item(14, 10).Value = 0.07552607652132566
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 12).Value = 0.6666666666666666
$ws.Cells.Item(14, 13).Value = 0.1474273333333333
$ws.Cells.Item(14, 14).Value = 0.442282
$ws.Cells.Item(14, 15).Value = 0.1588601259223368
$ws.Cells.Item(14, 16).Value = 0.1588601259223368
$ws.Cells.Item(14, 17).Value = 0.2053140369148889
$ws.Cells.Item(14, 18).Value = 1.847826332234
$ws.Cells.Item(14, 19).Value = 0.01199808202659783
$ws.Cells.Item(14, 20).Value = 0.01199808202659784
$ws.Cells.Item(15, 7).Value = 1.392645666666667
$ws.Cells.Item(15, 8).Value = 4.177937
$ws.Cells.Item(15, 9).Value = 0.07552607652132563
$ws.Cells.Item(15, 10).Value = 0.07552607652132566
$ws.Cells.Item(15, 15).Value = 0.4626735347223893
$ws.Cells.Item(15, 16).Value = 0.4626735347223893
$ws.Cells.Item(15, 17).Value = 0.5979686257706667
$ws.Cells.Item(15, 18).Value = 5.381717631936
$ws.Cells.Item(15, 19).Value = 0.03494391678783539
$ws.Cells.Item(15, 20).Value = 0.0349439167878354
$ws.Cells.Item(16, 7).Value = 1.392645666666667
$ws.Cells.Item(16, 8).Value = 4.177937
$ws.Cells.Item(16, 9).Value = 0.07552607652132563
$ws.Cells.Item(16, 10).Value = 0.07552607652132566
$ws.Cells.Item(16, 13).Value = 0.01780266666666666
$ws.Cells.Item(16, 14).Value = 0.053408
$ws.Cells.Item(16, 15).Value = 0.01918323966442261
$ws.Cells.Item(16, 16).Value = 0.01918323966442261
$ws.Cells.Item(16, 17).Value = 0.02479280658844444
$ws.Cells.Item(16, 18).Value = 0.223135259296
$ws.Cells.Item(16, 19).Value = 0.001448834826822111
$ws.Cells.Item(16, 20).Value = 0.001448834826822112
$ws.Cells.Item(17, 7).Value = 1.392645666666667
$ws.Cells.Item(17, 8).Value = 4.177937
$ws.Cells.Item(17, 9).Value = 0.07552607652132563
$ws.Cells.Item(17, 10).Value = 0.07552607652132566
$ws.Cells.Item(17, 13).Value = 0.3334263333333333
$ws.Cells.Item(17, 14).Value = 1.000279
$ws.Cells.Item(17, 15).Value = 0.3592830996908513
$ws.Cells.Item(17, 16).Value = 0.3592830996908513
$ws.Cells.Item(17, 17).Value = 0.4643447382692222
$ws.Cells.Item(17, 18).Value = 4.179102644423
$ws.Cells.Item(17, 19).Value = 0.0271352428800703
$ws.Cells.Item(17, 20).Value = 0.02713524288007031
$ws.Cells.Item(18, 7).Value = 1.139971333333333
$ws.Cells.Item(18, 8).Value = 3.419914
$ws.Cells.Item(18, 9).Value = 0.06182302089771886
$ws.Cells.Item(18, 10).Value = 0.06182302089771888
$ws.Cells.Item(18, 11).Value = 2
$ws.Cells.Item(18, 12).Value = 0.6666666666666666
$ws.Cells.Item(18, 13).Value = 0.1474273333333333
$ws.Cells.Item(18, 14).Value = 0.442282
$ws.Cells.Item(18, 15).Value = 0.1588601259223368
$ws.Cells.Item(18, 16).Value = 0.1588601259223368
$ws.Cells.Item(18, 17).Value = 0.1680629337497778
$ws.Cells.Item(18, 18).Value = 1.512566403748
$ws.Cells.Item(18, 19).Value = 0.009821212884710875
$ws.Cells.Item(18, 20).Value = 0.009821212884710879
$ws.Cells.Item(19, 7).Value = 1.139971333333333
$ws.Cells.Item(19, 8).Value = 3.419914
$ws.Cells.Item(19, 9).Value = 0.06182302089771886
$ws.Cells.Item(19, 10).Value = 0.06182302089771888
$ws.Cells.Item(19, 15).Value = 0.4626735347223893
$ws.Cells.Item(19, 16).Value = 0.4626735347223893
$ws.Cells.Item(19, 17).Value = 0.4894763312213334
$ws.Cells.Item(19, 18).Value = 4.405286980992
$ws.Cells.Item(19, 19).Value = 0.02860387560596373
$ws.Cells.Item(19, 20).Value = 0.02860387560596374
$ws.Cells.Item(20, 7).Value = 1.139971333333333
$ws.Cells.Item(20, 8).Value = 3.419914
$ws.Cells.Item(20, 9).Value = 0.06182302089771886
$ws.Cells.Item(20, 10).Value = 0.06182302089771888
$ws.Cells.Item(20, 13).Value = 0.01780266666666666
$ws.Cells.Item(20, 14).Value = 0.053408
$ws.Cells.Item(20, 15).Value = 0.01918323966442261
$ws.Cells.Item(20, 16).Value = 0.01918323966442261
$ws.Cells.Item(20, 17).Value = 0.02029452965688889
$ws.Cells.Item(20, 18).Value = 0.182650766912
$ws.Cells.Item(20, 19).Value = 0.001185965826659548
$ws.Cells.Item(20, 20).Value = 0.001185965826659549
$ws.Cells.Item(21, 7).Value = 1.139971333333333
$ws.Cells.Item(21, 8).Value = 3.419914
$ws.Cells.Item(21, 9).Value = 0.06182302089771886
$ws.Cells.Item(21, 10).Value = 0.06182302089771888
$ws.Cells.Item(21, 13).Value = 0.3334263333333333
$ws.Cells.Item(21, 14).Value = 1.000279
$ws.Cells.Item(21, 15).Value = 0.3592830996908513
$ws.Cells.Item(21, 16).Value = 0.3592830996908513
$ws.Cells.Item(21, 17).Value = 0.3800964617784445
$ws.Cells.Item(21, 18).Value = 3.420868156006
$ws.Cells.Item(21, 19).Value = 0.02221196658038471
$ws.Cells.Item(21, 20).Value = 0.02221196658038472
$ws.Cells.Item(22, 7).Value = 0.248809
$ws.Cells.Item(22, 8).Value = 0.746427
$ws.Cells.Item(22, 9).Value = 0.0134934305422948
$ws.Cells.Item(22, 10).Value = 0.01349343054229481
$ws.Cells.Item(22, 11).Value = 2
$ws.Cells.Item(22, 12).Value = 0.6666666666666666
$ws.Cells.Item(22, 13).Value = 0.1474273333333333
$ws.Cells.Item(22, 14).Value = 0.442282
$ws.Cells.Item(22, 15).Value = 0.1588601259223368
$ws.Cells.Item(22, 16).Value = 0.1588601259223368
$ws.Cells.Item(22, 17).Value = 0.03668124737933333
$ws.Cells.Item(22, 18).Value = 0.330131226414
$ws.Cells.Item(22, 19).Value = 0.002143568075073257
$ws.Cells.Item(22, 20).Value = 0.002143568075073258
$ws.Cells.Item(23, 7).Value = 0.248809
$ws.Cells.Item(23, 8).Value = 0.746427
$ws.Cells.Item(23, 9).Value = 0.0134934305422948
$ws.Cells.Item(23, 10).Value = 0.01349343054229481
$ws.Cells.Item(23, 15).Value = 0.4626735347223893
$ws.Cells.Item(23, 16).Value = 0.4626735347223893
$ws.Cells.Item(23, 17).Value = 0.106832613184
$ws.Cells.Item(23, 18).Value = 0.9614935186559999
$ws.Cells.Item(23, 19).Value = 0.006243053204534583
$ws.Cells.Item(23, 20).Value = 0.006243053204534585
$ws.Cells.Item(24, 7).Value = 0.248809
$ws.Cells.Item(24, 8).Value = 0.746427
$ws.Cells.Item(24, 9).Value = 0.0134934305422948
$ws.Cells.Item(24, 10).Value = 0.01349343054229481
$ws.Cells.Item(24, 13).Value = 0.01780266666666666
$ws.Cells.Item(24, 14).Value = 0.053408
$ws.Cells.Item(24, 15).Value = 0.01918323966442261
$ws.Cells.Item(24, 16).Value = 0.01918323966442261
$ws.Cells.Item(24, 17).Value = 0.004429463690666666
$ws.Cells.Item(24, 18).Value = 0.03986517321599999
$ws.Cells.Item(24, 19).Value = 0.0002588477119880812
$ws.Cells.Item(24, 20).Value = 0.0002588477119880813
$ws.Cells.Item(25, 7).Value = 0.248809
$ws.Cells.Item(25, 8).Value = 0.746427
$ws.Cells.Item(25, 9).Value = 0.0134934305422948
$ws.Cells.Item(25, 10).Value = 0.01349343054229481
$ws.Cells.Item(25, 13).Value = 0.3334263333333333
$ws.Cells.Item(25, 14).Value = 1.000279
$ws.Cells.Item(25, 15).Value = 0.3592830996908513
$ws.Cells.Item(25, 16).Value = 0.3592830996908513
$ws.Cells.Item(25, 17).Value = 0.08295947257033333
$ws.Cells.Item(25, 18).Value = 0.7466352531329999
$ws.Cells.Item(25, 19).Value = 0.004847961550698881
$ws.Cells.Item(25, 20).Value = 0.004847961550698882
